# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.565.02'
$ws.Range("E2").Value = '  +3.75%  '
$ws.Range("D3").Value = '2.426.01'
$ws.Range("E3").Value = '  +2.62%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '313.13'
$ws.Range("E6").Value = '  +6.23%  '
$ws.Range("E7").Value = '  +2.02%  '
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  +5.13%  '
$ws.Range("D10").Value = '35.24'
$ws.Range("E10").Value = '  +3.69%  '
$ws.Range("E11").Value = '  +2.26%  '
$ws.Range("E12").Value = '  +1.87%  '
$ws.Range("D13").Value = '18.87'
$ws.Range("E13").Value = '  +2.63%  '
$ws.Range("E14").Value = '  +3.15%  '
$ws.Range("D15").Value = '2.804.70'
$ws.Range("E15").Value = '  +2.66%  '
$ws.Range("D16").Value = '2.438.01'
$ws.Range("E16").Value = '  +3.75%  '
$ws.Range("D17").Value = '0.834'
$ws.Range("E17").Value = '  +5.05%  '
$ws.Range("D18").Value = '44.443.51'
$ws.Range("E18").Value = '  +3.49%  '
$ws.Range("D19").Value = '12.51'
$ws.Range("E19").Value = '  +4.79%  '
$ws.Range("E20").Value = '  +2.29%  '
$ws.Range("D21").Value = '0.0₃0934'
$ws.Range("E21").Value = '  +5.63%  '
$ws.Range("D22").Value = '68.96'
$ws.Range("E22").Value = '  +1.52%  '
$ws.Range("D23").Value = '241.00'
$ws.Range("E23").Value = '  +2.59%  '
$ws.Range("E24").Value = '  +5.13%  '
$ws.Range("D25").Value = '2.47'
$ws.Range("E25").Value = '  +1.64%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = '25.19'
$ws.Range("E27").Value = '  +2.49%  '
$ws.Range("E28").Value = '  -4.27%  '
$ws.Range("E29").Value = '  +3.92%  '
$ws.Range("D30").Value = '33.29'
$ws.Range("E30").Value = '  +5.60%  '
$ws.Range("D31").Value = '48.53'
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("D32").Value = '0.123'
$ws.Range("E32").Value = '  +18.29%  '
$ws.Range("D33").Value = '19.55'
$ws.Range("E33").Value = '  +12.94%  '
$ws.Range("E34").Value = '  +3.40%  '
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").Value = '0.0768'
$ws.Range("E36").Value = '  +6.97%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.54'
$ws.Range("E37").Value = '  +4.24%  '
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").Value = '1.89'
$ws.Range("E38").Value = '  +2.90%  '
$ws.Range("E39").Value = '  +4.39%  '
$ws.Range("D40").Value = '124.17'
$ws.Range("E40").Value = '  -1.82%  '
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("E42").Value = '  -4.14%  '
$ws.Range("D43").Value = '21.38'
$ws.Range("E43").Value = '  -0.72%  '
$ws.Range("E44").Value = '  +3.97%  '
$ws.Range("D45").Value = '1.949.95'
$ws.Range("E45").Value = '  +0.81%  '
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("D47").Value = '2.94'
$ws.Range("E47").Value = '  +8.52%  '
$ws.Range("D48").Value = '9.54'
$ws.Range("E48").Value = '  +3.62%  '
$ws.Range("E49").Value = '  +9.85%  '
$ws.Range("D50").Value = '2.669.86'
$ws.Range("E50").Value = '  +3.08%  '
$ws.Range("D51").Value = '54.13'
$ws.Range("E51").Value = '  +5.00%  '
